$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tiles")

# Clear any existing formulas in A2:A4 and set plain values for A1:B4
$ws.Range("A1").Value = 5
$ws.Range("B1").Value = 3

$ws.Range("A2").Formula = $null
$ws.Range("A2").Value = 15
$ws.Range("B2").Value = 2

$ws.Range("A3").Formula = $null
$ws.Range("A3").Value = 25
$ws.Range("B3").Value = 2

$ws.Range("A4").Formula = $null
$ws.Range("A4").Value = 35
$ws.Range("B4").Value = 2

# Update the selection to match the target (activeCell B5, sqref B5)
$ws.Range("B5").Select()
